# Generate Report for Handback
# Adds a new handback entry (file af7491df-2eb1-4a16-ae9a-4524175fb88b.md) as
# row 4 on the "Overview", "zh-cn" and "de-de" sheets, mirroring the existing
# rows 2/3 layout (Source File Name / File Extension / Status / Correspond
# Handoff File+Datetime / Target File / Correspond Handback File+DateTime /
# Handoff Reason), each with the matching hyperlink.

$wb = $excel.ActiveWorkbook

$guid = "af7491df-2eb1-4a16-ae9a-4524175fb88b"
$xlfHash = "32afd2df07e627058ae77c01644be8fc466ff6c5"

$statusInSync = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------
# Sheet "Overview" -> new row 4
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Hyperlinks.Add(
    $wsOverview.Cells.Item(4, 1),
    "https://github.com/OpenLocalizationTest/oltest/blob/9975c97dc5a059b55f03d084f517829aac13607b/e2e/$guid.md",
    [Type]::Missing,
    [Type]::Missing,
    "$guid.md")

$wsOverview.Cells.Item(4, 2).Value = $statusInSync
$wsOverview.Cells.Item(4, 3).Value = $statusInSync

# ---------------------------------------------------------------------
# Sheet "zh-cn" -> new row 4
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Cells.Item(4, 1),
    "https://github.com/OpenLocalizationTest/oltest/blob/9975c97dc5a059b55f03d084f517829aac13607b/e2e/$guid.md",
    [Type]::Missing,
    [Type]::Missing,
    "$guid.md")

$wsZhCn.Cells.Item(4, 2).Value = ".md"
$wsZhCn.Cells.Item(4, 3).Value = $statusInSync

$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Cells.Item(4, 4),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/f3a23ba4b4accfb3afb71c5340443158ef0ea1e4/ol-handoff/OpenLocalizationTestOrg/oltest-zhcn-fly/yuwzho/ht/$guid.$xlfHash.zh-cn.xlf",
    [Type]::Missing,
    [Type]::Missing,
    "$guid.$xlfHash.zh-cn.xlf")

$wsZhCn.Cells.Item(4, 5).Value = "2016-03-25 09:46:28"
$wsZhCn.Cells.Item(4, 5).NumberFormat = "yyyy-mm-dd HH:mm:ss"

$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Cells.Item(4, 6),
    "https://github.com/OpenLocalizationTestOrg/oltest-zhcn-fly/blob/edfeaa5ef65416ca356525118f45714436810cdd/e2e/$guid.md",
    [Type]::Missing,
    [Type]::Missing,
    "$guid.md")

$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Cells.Item(4, 7),
    "https://github.com/OpenLocalizationTestOrg/olhandback/blob/43eebafd80de212686aee9063c98aa85896e9657/ol-handback/OpenLocalizationTestOrg/oltest-zhcn-fly/yuwzho/ht/$guid.$xlfHash.zh-cn.xlf",
    [Type]::Missing,
    [Type]::Missing,
    "$guid.$xlfHash.zh-cn.xlf")

$wsZhCn.Cells.Item(4, 8).Value = "2016-03-25 09:47:17"
$wsZhCn.Cells.Item(4, 8).NumberFormat = "yyyy-mm-dd HH:mm:ss"

$wsZhCn.Cells.Item(4, 10).Value = "Include"

# ---------------------------------------------------------------------
# Sheet "de-de" -> new row 4
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Cells.Item(4, 1),
    "https://github.com/OpenLocalizationTest/oltest/blob/9975c97dc5a059b55f03d084f517829aac13607b/e2e/$guid.md",
    [Type]::Missing,
    [Type]::Missing,
    "$guid.md")

$wsDeDe.Cells.Item(4, 2).Value = ".md"
$wsDeDe.Cells.Item(4, 3).Value = $statusInSync

$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Cells.Item(4, 4),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/596d8218a363f3badd24bc3b34bd2af5b907aea7/ol-handoff/OpenLocalizationTestOrg/oltest-dede-fly/yuwzho/ht/$guid.$xlfHash.de-de.xlf",
    [Type]::Missing,
    [Type]::Missing,
    "$guid.$xlfHash.de-de.xlf")

$wsDeDe.Cells.Item(4, 5).Value = "2016-03-25 09:46:39"
$wsDeDe.Cells.Item(4, 5).NumberFormat = "yyyy-mm-dd HH:mm:ss"

$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Cells.Item(4, 6),
    "https://github.com/OpenLocalizationTestOrg/oltest-dede-fly/blob/ad11b14adaf25486bdc0e8ee8b34e86029f822ad/e2e/$guid.md",
    [Type]::Missing,
    [Type]::Missing,
    "$guid.md")

$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Cells.Item(4, 7),
    "https://github.com/OpenLocalizationTestOrg/olhandback/blob/c3b2df9728f42b35ebaf70a4ea6326389c653cd3/ol-handback/OpenLocalizationTestOrg/oltest-dede-fly/yuwzho/ht/$guid.$xlfHash.de-de.xlf",
    [Type]::Missing,
    [Type]::Missing,
    "$guid.$xlfHash.de-de.xlf")

$wsDeDe.Cells.Item(4, 8).Value = "2016-03-25 09:47:31"
$wsDeDe.Cells.Item(4, 8).NumberFormat = "yyyy-mm-dd HH:mm:ss"

$wsDeDe.Cells.Item(4, 10).Value = "Include"

Write-Output "Handback report row added for $guid"
